# Insert a new weekly price record for "Haba" (Vega Monumental Concepción)
# as row 55, pushing the existing rows 55-67 down to 56-68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 55 (shifts rows 55..67 -> 56..68)
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new record's values
$ws.Cells.Item(55, 1).Value  = 11
$ws.Cells.Item(55, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(55, 3).Value  = "Bíobío"
$ws.Cells.Item(55, 4).Value  = 45204
$ws.Cells.Item(55, 5).Value  = 8
$ws.Cells.Item(55, 6).Value  = 100112026
$ws.Cells.Item(55, 7).Value  = "Haba"
$ws.Cells.Item(55, 8).Value  = "Sin especificar"
$ws.Cells.Item(55, 9).Value  = "Primera"
$ws.Cells.Item(55, 10).Value = 50
$ws.Cells.Item(55, 11).Value = 14000
$ws.Cells.Item(55, 12).Value = 14000
$ws.Cells.Item(55, 13).Value = 14000
$ws.Cells.Item(55, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(55, 15).Value = "Región Metropolitana"
$ws.Cells.Item(55, 16).Value = 560
$ws.Cells.Item(55, 17).Value = 25
$ws.Cells.Item(55, 18).Value = "Hortaliza"
